# Add packages/footprints to BOM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column B: title / metadata block (values only shift within column B,
# the D:H table columns keep their existing row numbers) -------------------
$ws.Range("B2").Value  = "SinESC Multi v2.3C" + [char]0x2014 + "Bill of Materials (BOM)"
$ws.Range("B3").Value  = "Source: SinESC\Multi Edition\SinESC-Multi\SinESC-Multi.sch"
$ws.Range("B4").Value  = "Generated on: 12/30/2020"
$ws.Range("B10").Value = "BOM Consolidation Ratio (Total:Unique): 3.261:1"
$ws.Range("B5").Value  = "Last Updated: 1/18/2021"
$ws.Range("B6").Value  = "Tool: Eeschema (5.1.6)-1"
$ws.Range("B7").Value  = "Generator: bom_csv_grouped_by_value.py"
$ws.Range("B8").Value  = "Total Component Count: 75"
$ws.Range("B9").Value  = "Unique Component Count: 23"

# --- New column I: Package/Footprint ---------------------------------------
$ws.Range("I2").Value  = "Package/Footprint"
$ws.Range("I3").Value  = "C0201"
$ws.Range("I4").Value  = "C0201"
$ws.Range("I5").Value  = "C0402"
$ws.Range("I6").Value  = "C0603"
$ws.Range("I7").Value  = "C0603"
$ws.Range("I8").Value  = "C0402"
$ws.Range("I9").Value  = "SOD-523"
$ws.Range("I10").Value = "L1210"
$ws.Range("I11").Value = "TDSON-8-1"
$ws.Range("I12").Value = "R0201"
$ws.Range("I13").Value = "R0201"
$ws.Range("I14").Value = "R1206"
$ws.Range("I15").Value = "R0201"
$ws.Range("I16").Value = "R0201"
$ws.Range("I17").Value = "R0201"
$ws.Range("I18").Value = "R0201"
$ws.Range("I19").Value = "R0201"
$ws.Range("I20").Value = "R0201"
$ws.Range("I21").Value = "LQFP48-7x7mm-P0.5mm"
$ws.Range("I22").Value = "SOT23-6"
$ws.Range("I23").Value = "QFN37-EP-7x7mm"
$ws.Range("I24").Value = "SOT23-6"
$ws.Range("I25").Value = "Resonator_SMD_muRata_CSTxExxV-3Pin_3.0x1.1mm"

# Bold header row formatting for the new table header cells.
$ws.Range("D2:I2").Font.Bold = $true

# Match the saved view/selection of the target workbook.
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("I20").Select()
